$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=12; DAMSLTag='ba'; DialogAct='Appreciation'}
    @{Row=23; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=26; DAMSLTag='%'; DialogAct='Uninterpretable'}
    @{Row=35; DAMSLTag='%'; DialogAct='Uninterpretable'}
    @{Row=40; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=41; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=48; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
    @{Row=71; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
    @{Row=77; DAMSLTag='aa'; DialogAct='Agree/Accept'}
    @{Row=82; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=112; DAMSLTag='%'; DialogAct='Uninterpretable'}
    @{Row=125; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
    @{Row=130; DAMSLTag='aa'; DialogAct='Agree/Accept'}
    @{Row=136; DAMSLTag='ba'; DialogAct='Appreciation'}
    @{Row=142; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'}
    @{Row=148; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=158; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
    @{Row=159; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
    @{Row=161; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=170; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
    @{Row=182; DAMSLTag='aa'; DialogAct='Agree/Accept'}
    @{Row=185; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
    @{Row=191; DAMSLTag='%'; DialogAct='Uninterpretable'}
    @{Row=201; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=208; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=212; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=213; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
    @{Row=217; DAMSLTag='b'; DialogAct='Acknowledge (Backchannel)'}
    @{Row=223; DAMSLTag='aa'; DialogAct='Agree/Accept'}
    @{Row=228; DAMSLTag='aa'; DialogAct='Agree/Accept'}
    @{Row=230; DAMSLTag='ba'; DialogAct='Appreciation'}
    @{Row=239; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
    @{Row=240; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=241; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=243; DAMSLTag='sv'; DialogAct='Statement-opinion'}
    @{Row=245; DAMSLTag='sd'; DialogAct='Statement-non-opinion'}
    @{Row=250; DAMSLTag='aa'; DialogAct='Agree/Accept'}
    @{Row=252; DAMSLTag='sv'; DialogAct='Statement-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}
